# "update db for mysql"
# Adds more specific "JSON Object of ..." request/response descriptions
# to the Tenant / Operator / Promotion / PromotionBalance API rows, and
# moves the active cell selection from D15 to F22.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Tenant rows (2-5) ---
$ws.Range("E2").Value = "JSON Object of tenant"
$ws.Range("E3").Value = "JSON Object of tenant"
$ws.Range("F4").Value = "JSON Object of tenant"
$ws.Range("F5").Value = "JSON Objects of tenant"

# --- Operator rows (6-9) ---
$ws.Range("E6").Value = "JSON Object of Operator"
$ws.Range("E7").Value = "JSON Object of Operator"
$ws.Range("F8").Value = "JSON Object of Operator"
$ws.Range("F9").Value = "JSON Object of Operator"

# --- Promotion rows (10-13) ---
$ws.Range("E10").Value = "JSON Object of Promotion"
$ws.Range("E11").Value = "JSON Object of Promotion"
$ws.Range("F12").Value = "JSON Object of Promotion"
$ws.Range("F13").Value = "JSON Object of Promotion"

# --- PromotionBalance row (14) ---
$ws.Range("F14").Value = "JSON Object of PromotionBalance"

# Update the saved selection to match the author's cursor position.
[void]$ws.Range("F22").Select()
